$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cl = $sm.CustomLayouts.Item(1)
Write-Output ($cl | Get-Member | Out-String)
